# Update stats for 2025-10 (row 23 of Sheet1)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B23").Value = 6343
$ws.Range("D23").Value = 5902823
$ws.Range("E23").Value = 930.6042881917074
$ws.Range("F23").Value = 8.836650652024701
$ws.Range("H23").Value = 26.48238782252852
